$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 12653.542
$ws.Range("J87").Value = 12962.174
$ws.Range("L87").Value = 12962.174
$ws.Range("N87").Value = -15458.174

$ws.Range("H90").Value = 12653.542
$ws.Range("J90").Value = 12962.174
$ws.Range("L90").Value = 38886.522
$ws.Range("N90").Value = -51366.522

$ws.Range("H116").Value = 16670585
$ws.Range("I116").Value = 3197.5557
$ws.Range("J116").Value = 41671668
$ws.Range("K116").Value = 3197.5557
$ws.Range("L116").Value = 41671668
$ws.Range("M116").Value = 244.4443000000001
$ws.Range("N116").Value = -41678552

$ws.Range("H135").Value = 3044.12
$ws.Range("I135").Value = 213.47368
$ws.Range("J135").Value = 12007.833
$ws.Range("K135").Value = 1921.26312
$ws.Range("L135").Value = 108070.497
$ws.Range("M135").Value = 613.7368799999999
$ws.Range("N135").Value = -113140.497

$ws.Range("H137").Value = 3340.7192
$ws.Range("I137").Value = 2849.1143
$ws.Range("J137").Value = 4122.8184
$ws.Range("K137").Value = 8547.3429
$ws.Range("L137").Value = 12368.4552
$ws.Range("M137").Value = -5997.3429
$ws.Range("N137").Value = -17468.4552

$ws.Range("H138").Value = 5121.8594
$ws.Range("I138").Value = 3259.0454
$ws.Range("J138").Value = 6097.619
$ws.Range("K138").Value = 9777.136200000001
$ws.Range("L138").Value = 18292.857
$ws.Range("M138").Value = -4637.136200000001
$ws.Range("N138").Value = -28572.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 100000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 100000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 100000
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -101466

$ws.Range("H61").Value = 3892.889
$ws.Range("I61").Value = 2718.4285
$ws.Range("J61").Value = 8003.5
$ws.Range("K61").Value = 2718.4285
$ws.Range("L61").Value = 8003.5
$ws.Range("M61").Value = -2506.4285
$ws.Range("N61").Value = -8427.5

$ws.Range("H102").Value = 2903.913
$ws.Range("I102").Value = 2312
$ws.Range("J102").Value = 4013.75
$ws.Range("K102").Value = 2312
$ws.Range("L102").Value = 4013.75
$ws.Range("M102").Value = -690
$ws.Range("N102").Value = -7257.75

$ws.Range("H136").Value = 3892.889
$ws.Range("I136").Value = 2718.4285
$ws.Range("J136").Value = 8003.5
$ws.Range("K136").Value = 8155.2855
$ws.Range("L136").Value = 24010.5
$ws.Range("M136").Value = -5605.2855
$ws.Range("N136").Value = -29110.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1598.2
$ws.Range("I105").Value = 1528.8889
$ws.Range("K105").Value = 1528.8889
$ws.Range("M105").Value = 218.1111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 16365.5
$ws.Range("I33").Value = 2731
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 2731
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -2352
$ws.Range("N33").Value = -30758

$ws.Range("H36").Value = 5349.6665
$ws.Range("I36").Value = 3219.6
$ws.Range("K36").Value = 3219.6
$ws.Range("M36").Value = -2831.6

$ws.Range("H40").Value = 5349.6665
$ws.Range("I40").Value = 3219.6
$ws.Range("K40").Value = 3219.6
$ws.Range("M40").Value = -3059.6

$ws.Range("H134").Value = 1161.7273
$ws.Range("I134").Value = 1076.3103
$ws.Range("J134").Value = 1781
$ws.Range("K134").Value = 3228.9309
$ws.Range("L134").Value = 5343
$ws.Range("M134").Value = -693.9309000000003
$ws.Range("N134").Value = -10413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1215.1538
$ws.Range("J34").Value = 1310.7
$ws.Range("L34").Value = 3932.1
$ws.Range("N34").Value = -4100.1

$ws.Range("H39").Value = 1042.4348
$ws.Range("J39").Value = 1058.8
$ws.Range("L39").Value = 3176.4
$ws.Range("N39").Value = -3764.4

$ws.Range("H55").Value = 2426.36
$ws.Range("J55").Value = 2476.3044
$ws.Range("L55").Value = 7428.9132
$ws.Range("N55").Value = -7782.9132

$ws.Range("H137").Value = 7579755.5
$ws.Range("I137").Value = 2297.1667
$ws.Range("J137").Value = 10421302
$ws.Range("K137").Value = 6891.500100000001
$ws.Range("L137").Value = 31263906
$ws.Range("M137").Value = -1791.500100000001
$ws.Range("N137").Value = -31274106

$ws.Range("H140").Value = 1920.7407
$ws.Range("I140").Value = 1124.2106
$ws.Range("K140").Value = 3372.6318
$ws.Range("M140").Value = 1807.3682

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3664.2856
$ws.Range("I122").Value = 2275.818
$ws.Range("J122").Value = 5191.6
$ws.Range("K122").Value = 6827.454000000001
$ws.Range("L122").Value = 15574.8
$ws.Range("M122").Value = -4377.454000000001
$ws.Range("N122").Value = -20474.8

$ws.Range("H132").Value = 1962.0667
$ws.Range("I132").Value = 1446.6
$ws.Range("J132").Value = 2993
$ws.Range("K132").Value = 4339.799999999999
$ws.Range("L132").Value = 8979
$ws.Range("M132").Value = -1809.799999999999
$ws.Range("N132").Value = -14039

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 452.34784
$ws.Range("I16").Value = 468.27274
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = 468.27274
$ws.Range("L16").Value = 102
$ws.Range("M16").Value = -298.27274
$ws.Range("N16").Value = -442

$ws.Range("H22").Value = 1048.5
$ws.Range("I22").Value = 764.6667
$ws.Range("J22").Value = 1900
$ws.Range("K22").Value = 764.6667
$ws.Range("L22").Value = 1900
$ws.Range("M22").Value = -469.6667
$ws.Range("N22").Value = -2490

$ws.Range("H27").Value = 1048.5
$ws.Range("I27").Value = 764.6667
$ws.Range("J27").Value = 1900
$ws.Range("K27").Value = 764.6667
$ws.Range("L27").Value = 1900
$ws.Range("M27").Value = -657.6667
$ws.Range("N27").Value = -2114

$ws.Range("H40").Value = 3536.0527
$ws.Range("I40").Value = 3152.3333
$ws.Range("J40").Value = 4975
$ws.Range("K40").Value = 3152.3333
$ws.Range("L40").Value = 4975
$ws.Range("M40").Value = -3016.3333
$ws.Range("N40").Value = -5247

$ws.Range("H55").Value = 470.79166
$ws.Range("I55").Value = 439.8
$ws.Range("J55").Value = 522.44446
$ws.Range("K55").Value = 439.8
$ws.Range("L55").Value = 522.44446
$ws.Range("M55").Value = -266.8
$ws.Range("N55").Value = -868.44446

$ws.Range("H68").Value = 2232.4407
$ws.Range("I68").Value = 995.3333
$ws.Range("J68").Value = 2775.561
$ws.Range("K68").Value = 995.3333
$ws.Range("L68").Value = 2775.561
$ws.Range("M68").Value = -246.3333
$ws.Range("N68").Value = -4273.561

$ws.Range("H71").Value = 2232.4407
$ws.Range("I71").Value = 995.3333
$ws.Range("J71").Value = 2775.561
$ws.Range("K71").Value = 4976.6665
$ws.Range("L71").Value = 13877.805
$ws.Range("M71").Value = -1232.6665
$ws.Range("N71").Value = -21365.805

$ws.Range("H132").Value = 5063.085
$ws.Range("I132").Value = 2475.2354
$ws.Range("J132").Value = 11831.308
$ws.Range("K132").Value = 7425.706200000001
$ws.Range("L132").Value = 35493.924
$ws.Range("M132").Value = -4895.706200000001
$ws.Range("N132").Value = -40553.924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 854.807
$ws.Range("I132").Value = 702.3774
$ws.Range("K132").Value = 2107.1322
$ws.Range("M132").Value = 422.8678

$ws.Range("H136").Value = 841.96155
$ws.Range("I136").Value = 821.34784
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2464.04352
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = 85.95647999999983
$ws.Range("N136").Value = -8100
